$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLOYEE DTR")

# Rows 5 through 10: fill in the Official Business time columns
# K = OFFICIAL BUSINESS DEPARTURE, L = OFFICIAL BUSINESS TIME START -> 08:30:00
# M = OFFICIAL BUSINESS TIME END, N = OFFICIAL BUSINESS ARRIVAL -> 18:30:00
for ($row = 5; $row -le 10; $row++) {
    $ws.Range("K$row").Value = "08:30:00"
    $ws.Range("L$row").Value = "08:30:00"
    $ws.Range("M$row").Value = "18:30:00"
    $ws.Range("N$row").Value = "18:30:00"
}
